$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

$aboutWs = $wb.Worksheets.Item("About")
$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")

# Update "About" sheet
$aboutWs.Range("A2").Value2 = "Version: " + $newVersion

$aboutWs.Range("A6").Value2 = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Leer Coal Mine, United States, M1036, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# Update "Boundaries and methane sources" sheet, column S, rows 2-11
for ($row = 2; $row -le 11; $row++) {
    $cell = $dataWs.Cells.Item($row, 19)
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value2 = $newVersion
    }
}

$wb.Save()
